# Generate Report for Handback
# Updates handback/handoff timestamps and priority values produced by a
# fresh report-generation run. Rows 4 and 5 on the language sheets (and
# the corresponding rows on Overview) shared identical text values, so
# updating the value also updates both rows that held that same text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" ---
$wsOverview.Range("G4").Value = "2016-08-13 14:18:39"
$wsOverview.Range("G5").Value = "2016-08-13 14:18:39"

# --- zh-cn sheet ---
$wsZhCn.Range("E4").Value = "mt"                        # Priority
$wsZhCn.Range("E5").Value = "mt"                        # Priority
$wsZhCn.Range("H4").Value = "2016-08-13 14:18:31"        # Correspond Handoff Datetime
$wsZhCn.Range("H5").Value = "2016-08-13 14:18:31"        # Correspond Handoff Datetime
$wsZhCn.Range("K4").Value = "2016-08-13 14:18:59"        # Correspond Handback DateTime
$wsZhCn.Range("K5").Value = "2016-08-13 14:18:59"        # Correspond Handback DateTime

# --- de-de sheet ---
$wsDeDe.Range("E4").Value = "mt"                        # Priority
$wsDeDe.Range("E5").Value = "mt"                        # Priority
$wsDeDe.Range("H4").Value = "2016-08-13 14:18:39"        # Correspond Handoff Datetime
$wsDeDe.Range("H5").Value = "2016-08-13 14:18:39"        # Correspond Handoff Datetime
$wsDeDe.Range("K4").Value = "2016-08-13 14:19:11"        # Correspond Handback DateTime
$wsDeDe.Range("K5").Value = "2016-08-13 14:19:11"        # Correspond Handback DateTime
